$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'307.96"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "'0.89%"
$ws.Range("E2").ClearFormats()
$ws.Range("D3").Value = "'36.35"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "'1.32%"
$ws.Range("E3").ClearFormats()
$ws.Range("D4").Value = "'5.061"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "'1.33%"
$ws.Range("E4").ClearFormats()
$ws.Range("D5").Value = "'0.08098"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "'0.36%"
$ws.Range("E5").ClearFormats()
$ws.Range("D6").Value = "'1.994"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "'5.49%"
$ws.Range("E6").ClearFormats()
$ws.Range("D7").Value = "'7.829"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "'-0.36%"
$ws.Range("E7").ClearFormats()
$ws.Range("D8").Value = "'0.9275"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "'-0.75%"
$ws.Range("E8").ClearFormats()
$ws.Range("D9").Value = "'0.1485"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "'12.67%"
$ws.Range("E9").ClearFormats()
$ws.Range("D10").Value = "'0.1936"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "'1.70%"
$ws.Range("E10").ClearFormats()
$ws.Range("D11").Value = "'0.09099"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "'-1.15%"
$ws.Range("E11").ClearFormats()
$ws.Range("D12").Value = "'0.03528"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "'0.39%"
$ws.Range("E12").ClearFormats()
$ws.Range("D13").Value = "'0.09866"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "'-0.41%"
$ws.Range("E13").ClearFormats()
$ws.Range("D14").Value = "'0.001417"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "'-0.62%"
$ws.Range("E14").ClearFormats()
$ws.Range("D15").Value = "'0.006316"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "'-5.93%"
$ws.Range("E15").ClearFormats()
$ws.Range("D16").Value = "'3.850"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "'6.91%"
$ws.Range("E16").ClearFormats()
$ws.Range("D17").Value = "'4.166"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "'0.31%"
$ws.Range("E17").ClearFormats()
$ws.Range("D19").Value = "'0.3449"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "'-0.02%"
$ws.Range("E19").ClearFormats()
$ws.Range("D20").Value = "'0.1328"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "'1.65%"
$ws.Range("E20").ClearFormats()
$ws.Range("D21").Value = "'4.805"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "'-8.44%"
$ws.Range("E21").ClearFormats()
$ws.Range("D22").Value = "'0.2348"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "'-7.34%"
$ws.Range("E22").ClearFormats()
$ws.Range("D23").Value = "'0.04395"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "'-0.59%"
$ws.Range("E23").ClearFormats()
$ws.Range("E24").Value = "'-0.02%"
$ws.Range("E24").ClearFormats()
$ws.Range("D25").Value = "'0.004156"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "'-11.97%"
$ws.Range("E25").ClearFormats()
$ws.Range("D27").Value = "'0.0001303"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "'0.11%"
$ws.Range("E27").ClearFormats()
$ws.Range("E39").Value = "'5.34%"
$ws.Range("E39").ClearFormats()
$ws.Range("D40").Value = "'0.05104"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "'-1.79%"
$ws.Range("E40").ClearFormats()
$ws.Range("D41").Value = "'0.007491"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "'-0.88%"
$ws.Range("E41").ClearFormats()
$ws.Range("D42").Value = "'0.01012"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "'-0.38%"
$ws.Range("E42").ClearFormats()
$ws.Range("E43").Value = "'-0.42%"
$ws.Range("E43").ClearFormats()
$ws.Range("D44").Value = "'0.002124"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "'-6.92%"
$ws.Range("E44").ClearFormats()
$ws.Range("D45").Value = "'0.009659"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "'-9.59%"
$ws.Range("E45").ClearFormats()
$ws.Range("D46").Value = "'0.00006303"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "'-0.39%"
$ws.Range("E46").ClearFormats()
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "'0.11%"
$ws.Range("E47").ClearFormats()
$ws.Range("E48").Value = "'-0.24%"
$ws.Range("E48").ClearFormats()
$ws.Range("D49").Value = "'0.001605"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "'-3.34%"
$ws.Range("E49").ClearFormats()
$ws.Range("E50").Value = "'0.11%"
$ws.Range("E50").ClearFormats()
$ws.Range("E51").Value = "'0.11%"
$ws.Range("E51").ClearFormats()
